$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 9949.7999999999993
$ws.Range("B5").Value = 10015.91
$ws.Range("C5").Value = 282.89999999999998
$ws.Range("D5").Value = 284.76
$ws.Range("E5").Value = $true
$ws.Range("F5").Value = 0.66
$ws.Range("G5").Value = 42609.506041666667
$ws.Range("G5").NumberFormat = "m/d/yy h:mm"
$ws.Range("H5").Value = $false
